# Weekly fruit/vegetable price update.
# A new weekly reading is inserted above the current row 28, which pushes
# the existing rows 28-33 down to 29-34 (row 33's data ends up at row 34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 28, shifting rows 28:33 down to 29:34.
$ws.Rows("28:28").Insert()

# Fill in the new row 28 with this week's reading.
$ws.Range("A28").Value = 1
$ws.Range("B28").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C28").Value = "Arica y Parinacota"
$ws.Range("D28").Value = 44943
$ws.Range("E28").Value = 15
$ws.Range("F28").Value = 100112003
$ws.Range("G28").Value = "Ajo"
$ws.Range("H28").Value = "Chino"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 400
$ws.Range("K28").Value = 16000
$ws.Range("L28").Value = 17000
$ws.Range("M28").Value = 16500
$ws.Range("N28").Value = "`$/caja 10 kilos"
$ws.Range("O28").Value = "China"
$ws.Range("P28").Value = 1650
$ws.Range("Q28").Value = 10
$ws.Range("R28").Value = "Hortaliza"
